$p = $ppt.ActivePresentation

# Update the remaining "Adaptive Warehouse Right-Sizing" slide: compress the
# four detailed bullets down to three concise ones.
$s1 = $p.Slides.Item(1)
$body = $s1.Shapes.Item(2)
$body.TextFrame.TextRange.Text = "Auto-rightsizes warehouses by hourly load`rPolicy DT derives size from credits_used`rExecutor task applies changes and logs"

# Remove the second "How it works" details slide entirely.
$p.Slides.Item(2).Delete()
